# Product Backlog update:
#  - reprioritize several existing backlog items
#  - add two new backlog rows (hosted pdf source, build project / grouping)
#  - sheet view zoom + selection changed, sheet data grid now spans to row 15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) down to the two new rows (14 and 15) by copying
# the format from the last existing data row before writing new values.
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D15").PasteSpecial(-4122)

# --- Row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "User"
$ws.Range("C2").Value = "Enter notes for a source"
$ws.Range("D2").Value = "I can remember key details or insights from the source"

# --- Row 3 ---
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "User"
$ws.Range("C3").Value = "Upload a pdf source"
$ws.Range("D3").Value = "I can enter notes for that source"

# --- Row 4 ---
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "User"
$ws.Range("C4").Value = "Upload a embeded video source"
$ws.Range("D4").Value = "I can enter notes for that source"

# --- Row 5 (unchanged content, rewritten for clarity/idempotency) ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "User"
$ws.Range("C5").Value = "Upload a hosted video source"
$ws.Range("D5").Value = "I can enter notes for that source"

# --- Row 6 (now the new "Upload a hosted pdf source" story) ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "User"
$ws.Range("C6").Value = "Upload a hosted pdf source"
$ws.Range("D6").Value = "I can enter notes for that source"

# --- Row 7 ---
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "User"
$ws.Range("C7").Value = "Enter a tag for a note"
$ws.Range("D7").Value = "I can organize my notes based on topic"

# --- Row 8 ---
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "User"
$ws.Range("C8").Value = "Find notes related to a specific tag"
$ws.Range("D8").Value = "I can find my notes related to that topic across all sources"

# --- Row 9 ---
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "User"
$ws.Range("C9").Value = "Access the system offline"
$ws.Range("D9").Value = "I can work when a network connection is not available (e.g., during plane flight)"

# --- Row 10 ---
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "User"
$ws.Range("C10").Value = "Share notes with another user"
$ws.Range("D10").Value = "I can share my work on a topic with another user for a shared project"

# --- Row 11 ---
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "User"
$ws.Range("C11").Value = "Add comments to a note"
$ws.Range("D11").Value = "I can followup on a previous note from either myself or another user"

# --- Row 12 ---
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "User"
$ws.Range("C12").Value = "Export a list of sources"
$ws.Range("D12").Value = "I can add the sources to an external document bibliography"

# --- Row 13 ---
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "User"
$ws.Range("C13").Value = "Attach a note to a specific portion of a source"
$ws.Range("D13").Value = "I can identify the relevant portion of the source for the specific note"

# --- Row 14 (new location of "Find sources related to a specific tag") ---
$ws.Range("A14").Value = 1.5
$ws.Range("B14").Value = "User"
$ws.Range("C14").Value = "Find sources related to a specific tag"
$ws.Range("D14").Value = "I can find a list of sources related to that topic"

# --- Row 15 (brand new "Build Project" story) ---
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "User"
$ws.Range("C15").Value = "Build Project"
$ws.Range("D15").Value = "I can group a set of sources while working on project"

# --- View state updates ---
$win = $excel.ActiveWindow
$win.Zoom = 113
$null = $ws.Range("D21").Select()

Write-Host "Product backlog updated"
